$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GRAFICO")

# --- New date pair (14/jun -> serial 43630) in the header row ---
$ws.Range("V3").Value = 43630
$ws.Range("W3").Value = 43630

# --- Re-format the previously "latest" pair (T4:U4 / T5:U5) so it matches
#     the plain "filled" look used by every earlier pair, by copying the
#     formatting already used for earlier filled pairs (H4:I4 / H5:I5). ---
$ws.Range("H4:I4").Copy()
$ws.Range("T4:U4").PasteSpecial(-4122)

$ws.Range("H5:I5").Copy()
$ws.Range("T5:U5").PasteSpecial(-4122)

# --- New pair of values becomes the new "latest" pair (V4:W4 / V5:W5); row 5
#     still needs the "filled" formatting applied (row 4's was already there). ---
$ws.Range("H5:I5").Copy()
$ws.Range("V5:W5").PasteSpecial(-4122)

$ws.Range("V4").Value = 4.03
$ws.Range("W4").Value = 1725.55

$ws.Range("V5").Value = 3.91
$ws.Range("W5").Value = 1757.42

$excel.CutCopyMode = $false

# --- Match the author's recorded view state for the sheet ---
$ws.Activate()
$excel.Goto($ws.Range("V13"), $true)
